$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 'Meal Inexpensive Restaurant'
$ws.Cells.Item(2, 3).Value = 1600
$ws.Cells.Item(2, 4).Value = 1100
$ws.Cells.Item(2, 5).Value = 2000

$ws.Cells.Item(3, 2).Value = 'Meal for 2 People Midrange Restaurant Threecourse'
$ws.Cells.Item(3, 3).Value = 10500
$ws.Cells.Item(3, 4).Value = 7000
$ws.Cells.Item(3, 5).Value = 12000

$ws.Cells.Item(4, 2).Value = 'McMeal at McDonalds or Equivalent Combo Meal'
$ws.Cells.Item(4, 3).Value = 1000
$ws.Cells.Item(4, 4).Value = 1000
$ws.Cells.Item(4, 5).Value = 1100

$ws.Cells.Item(5, 2).Value = 'Domestic Beer 05 liter draught'
$ws.Cells.Item(5, 3).Value = 900
$ws.Cells.Item(5, 4).Value = 600
$ws.Cells.Item(5, 5).Value = 1000

$ws.Cells.Item(6, 2).Value = 'Imported Beer 033 liter bottle'
$ws.Cells.Item(6, 3).Value = 900
$ws.Cells.Item(6, 4).Value = 600
$ws.Cells.Item(6, 5).Value = 1000

$ws.Cells.Item(7, 2).Value = 'Cappuccino regular'
$ws.Cells.Item(7, 3).Value = 479
$ws.Cells.Item(7, 4).Value = 350
$ws.Cells.Item(7, 5).Value = 600

$ws.Cells.Item(8, 2).Value = 'CokePepsi 033 liter bottle'
$ws.Cells.Item(8, 3).Value = 295
$ws.Cells.Item(8, 4).Value = 205
$ws.Cells.Item(8, 5).Value = 400

$ws.Cells.Item(9, 2).Value = 'Water 033 liter bottle'
$ws.Cells.Item(9, 3).Value = 243
$ws.Cells.Item(9, 4).Value = 200
$ws.Cells.Item(9, 5).Value = 300

$ws.Cells.Item(10, 2).Value = 'Milk regular 1 liter'
$ws.Cells.Item(10, 3).Value = 115
$ws.Cells.Item(10, 4).Value = 83
$ws.Cells.Item(10, 5).Value = 150

$ws.Cells.Item(11, 2).Value = 'Loaf of Fresh White Bread 500g'
$ws.Cells.Item(11, 3).Value = 350
$ws.Cells.Item(11, 4).Value = 139
$ws.Cells.Item(11, 5).Value = 500

$ws.Cells.Item(12, 2).Value = 'Rice white 1kg'
$ws.Cells.Item(12, 3).Value = 286
$ws.Cells.Item(12, 4).Value = 200
$ws.Cells.Item(12, 5).Value = 435

$ws.Cells.Item(13, 2).Value = 'Eggs regular 12'
$ws.Cells.Item(13, 3).Value = 353
$ws.Cells.Item(13, 4).Value = 203
$ws.Cells.Item(13, 5).Value = 480

$ws.Cells.Item(14, 2).Value = 'Local Cheese 1kg'
$ws.Cells.Item(14, 3).Value = 875
$ws.Cells.Item(14, 4).Value = 595
$ws.Cells.Item(14, 5).Value = 1400

$ws.Cells.Item(15, 2).Value = 'Chicken Fillets 1kg'
$ws.Cells.Item(15, 3).Value = 1218
$ws.Cells.Item(15, 4).Value = 799
$ws.Cells.Item(15, 5).Value = 1400

$ws.Cells.Item(16, 2).Value = 'Beef Round 1kg or Equivalent Back Leg Red Meat'
$ws.Cells.Item(16, 3).Value = 2005
$ws.Cells.Item(16, 4).Value = 1000
$ws.Cells.Item(16, 5).Value = 2500

$ws.Cells.Item(17, 2).Value = 'Apples 1kg'
$ws.Cells.Item(17, 3).Value = 292
$ws.Cells.Item(17, 4).Value = 100
$ws.Cells.Item(17, 5).Value = 400

$ws.Cells.Item(18, 2).Value = 'Banana 1kg'
$ws.Cells.Item(18, 3).Value = 185
$ws.Cells.Item(18, 4).Value = 125
$ws.Cells.Item(18, 5).Value = 500

$ws.Cells.Item(19, 2).Value = 'Oranges 1kg'
$ws.Cells.Item(19, 3).Value = 297
$ws.Cells.Item(19, 4).Value = 120
$ws.Cells.Item(19, 5).Value = 400

$ws.Cells.Item(20, 2).Value = 'Tomato 1kg'
$ws.Cells.Item(20, 3).Value = 457
$ws.Cells.Item(20, 4).Value = 220
$ws.Cells.Item(20, 5).Value = 600

$ws.Cells.Item(21, 2).Value = 'Potato 1kg'
$ws.Cells.Item(21, 3).Value = 129
$ws.Cells.Item(21, 4).Value = 75
$ws.Cells.Item(21, 5).Value = 240

$ws.Cells.Item(22, 2).Value = 'Onion 1kg'
$ws.Cells.Item(22, 3).Value = 172
$ws.Cells.Item(22, 4).Value = 130
$ws.Cells.Item(22, 5).Value = 271

$ws.Cells.Item(23, 2).Value = 'Lettuce 1 head'
$ws.Cells.Item(23, 3).Value = 188
$ws.Cells.Item(23, 4).Value = 139
$ws.Cells.Item(23, 5).Value = 220

$ws.Cells.Item(24, 2).Value = 'Water 15 liter bottle'
$ws.Cells.Item(24, 3).Value = 208
$ws.Cells.Item(24, 4).Value = 117
$ws.Cells.Item(24, 5).Value = 250

$ws.Cells.Item(25, 2).Value = 'Bottle of Wine MidRange'
$ws.Cells.Item(25, 3).Value = 1500
$ws.Cells.Item(25, 4).Value = 1000
$ws.Cells.Item(25, 5).Value = 1800

$ws.Cells.Item(26, 2).Value = 'Domestic Beer 05 liter bottle'
$ws.Cells.Item(26, 3).Value = 351
$ws.Cells.Item(26, 4).Value = 250
$ws.Cells.Item(26, 5).Value = 400

$ws.Cells.Item(27, 2).Value = 'Imported Beer 033 liter bottle'
$ws.Cells.Item(27, 3).Value = 414
$ws.Cells.Item(27, 4).Value = 300
$ws.Cells.Item(27, 5).Value = 500

$ws.Cells.Item(28, 2).Value = 'Cigarettes 20 Pack Marlboro'
$ws.Cells.Item(28, 3).Value = 1100
$ws.Cells.Item(28, 4).Value = 995
$ws.Cells.Item(28, 5).Value = 1100

$ws.Cells.Item(29, 2).Value = 'Oneway Ticket Local Transport'
$ws.Cells.Item(29, 3).Value = 300
$ws.Cells.Item(29, 4).Value = 280
$ws.Cells.Item(29, 5).Value = 310

$ws.Cells.Item(30, 2).Value = 'Monthly Pass Regular Price'
$ws.Cells.Item(30, 3).Value = 7000
$ws.Cells.Item(30, 4).Value = 6000
$ws.Cells.Item(30, 5).Value = 9900

$ws.Cells.Item(31, 2).Value = 'Taxi Start Normal Tariff'
$ws.Cells.Item(31, 3).Value = 800
$ws.Cells.Item(31, 4).Value = 650
$ws.Cells.Item(31, 5).Value = 1000

$ws.Cells.Item(32, 2).Value = 'Taxi 1km Normal Tariff'
$ws.Cells.Item(32, 3).Value = 150
$ws.Cells.Item(32, 4).Value = 120
$ws.Cells.Item(32, 5).Value = 199

$ws.Cells.Item(33, 2).Value = 'Taxi 1hour Waiting Normal Tariff'
$ws.Cells.Item(33, 3).Value = 4600
$ws.Cells.Item(33, 4).Value = 4440
$ws.Cells.Item(33, 5).Value = 5500

$ws.Cells.Item(34, 2).Value = 'Gasoline 1 liter'
$ws.Cells.Item(34, 3).Value = 194
$ws.Cells.Item(34, 4).Value = 181
$ws.Cells.Item(34, 5).Value = 220

$ws.Cells.Item(35, 2).Value = 'Volkswagen Golf 14 90 KW Trendline Or Equivalent New Car'
$ws.Cells.Item(35, 3).Value = 2850000
$ws.Cells.Item(35, 4).Value = 2800000
$ws.Cells.Item(35, 5).Value = 2855000

$ws.Cells.Item(36, 2).Value = 'Toyota Corolla Sedan 16l 97kW Comfort Or Equivalent New Car'
$ws.Cells.Item(36, 3).Value = 3021414
$ws.Cells.Item(36, 4).Value = 3000000
$ws.Cells.Item(36, 5).Value = 3000000

$ws.Cells.Item(37, 1).Value = 'Utilities Monthly'
$ws.Cells.Item(37, 2).Value = 'Basic Electricity Heating Cooling Water Garbage for 85m2 Apartment'
$ws.Cells.Item(37, 3).Value = 10601
$ws.Cells.Item(37, 4).Value = 7600
$ws.Cells.Item(37, 5).Value = 20000

$ws.Cells.Item(38, 1).Value = 'Utilities Monthly'
$ws.Cells.Item(38, 2).Value = 'Mobile Phone Monthly Plan with Calls and 10GB Data'
$ws.Cells.Item(38, 3).Value = 2685
$ws.Cells.Item(38, 4).Value = 2000
$ws.Cells.Item(38, 5).Value = 2699

$ws.Cells.Item(39, 1).Value = 'Utilities Monthly'
$ws.Cells.Item(39, 2).Value = 'Internet 60 Mbps or More Unlimited Data CableADSL'
$ws.Cells.Item(39, 3).Value = 1813
$ws.Cells.Item(39, 4).Value = 1000
$ws.Cells.Item(39, 5).Value = 2500

$ws.Cells.Item(40, 2).Value = 'Fitness Club Monthly Fee for 1 Adult'
$ws.Cells.Item(40, 3).Value = 5341
$ws.Cells.Item(40, 4).Value = 2900
$ws.Cells.Item(40, 5).Value = 6500

$ws.Cells.Item(41, 2).Value = 'Tennis Court Rent 1 Hour on Weekend'
$ws.Cells.Item(41, 3).Value = 3244
$ws.Cells.Item(41, 4).Value = 2500
$ws.Cells.Item(41, 5).Value = 3500

$ws.Cells.Item(42, 2).Value = 'Cinema International Release 1 Seat'
$ws.Cells.Item(42, 3).Value = 1600
$ws.Cells.Item(42, 4).Value = 1400
$ws.Cells.Item(42, 5).Value = 1800

$ws.Cells.Item(43, 2).Value = 'Preschool or Kindergarten Full Day Private Monthly for 1 Child'
$ws.Cells.Item(43, 3).Value = 45600
$ws.Cells.Item(43, 4).Value = 23000
$ws.Cells.Item(43, 5).Value = 103700

$ws.Cells.Item(44, 2).Value = 'International Primary School Yearly for 1 Child'
$ws.Cells.Item(44, 3).Value = 1825000
$ws.Cells.Item(44, 4).Value = 1600000
$ws.Cells.Item(44, 5).Value = 2400000

$ws.Cells.Item(45, 2).Value = '1 Pair of Jeans Levis 501 Or Similar'
$ws.Cells.Item(45, 3).Value = 9724
$ws.Cells.Item(45, 4).Value = 5500
$ws.Cells.Item(45, 5).Value = 11000

$ws.Cells.Item(46, 2).Value = '1 Summer Dress in a Chain Store Zara HM '
$ws.Cells.Item(46, 3).Value = 4876
$ws.Cells.Item(46, 4).Value = 4000
$ws.Cells.Item(46, 5).Value = 6500

$ws.Cells.Item(47, 2).Value = '1 Pair of Nike Running Shoes MidRange'
$ws.Cells.Item(47, 3).Value = 10235
$ws.Cells.Item(47, 4).Value = 7900
$ws.Cells.Item(47, 5).Value = 13000

$ws.Cells.Item(48, 3).Value = 14871
$ws.Cells.Item(48, 4).Value = 11000
$ws.Cells.Item(48, 5).Value = 16000

$ws.Cells.Item(49, 2).Value = 'Apartment 1 bedroom in City Centre'
$ws.Cells.Item(49, 3).Value = 106724
$ws.Cells.Item(49, 4).Value = 89900
$ws.Cells.Item(49, 5).Value = 125000

$ws.Cells.Item(50, 2).Value = 'Apartment 1 bedroom Outside of Centre'
$ws.Cells.Item(50, 3).Value = 85129
$ws.Cells.Item(50, 4).Value = 65000
$ws.Cells.Item(50, 5).Value = 100000

$ws.Cells.Item(51, 2).Value = 'Apartment 3 bedrooms in City Centre'
$ws.Cells.Item(51, 3).Value = 172226
$ws.Cells.Item(51, 4).Value = 150000
$ws.Cells.Item(51, 5).Value = 200000

$ws.Cells.Item(52, 2).Value = 'Apartment 3 bedrooms Outside of Centre'
$ws.Cells.Item(52, 3).Value = 141406
$ws.Cells.Item(52, 4).Value = 110000
$ws.Cells.Item(52, 5).Value = 165000

$ws.Cells.Item(53, 3).Value = 765788
$ws.Cells.Item(53, 4).Value = 600000
$ws.Cells.Item(53, 5).Value = 900000

$ws.Cells.Item(54, 3).Value = 516325
$ws.Cells.Item(54, 4).Value = 350000
$ws.Cells.Item(54, 5).Value = 600000

$ws.Cells.Item(55, 2).Value = 'Average Monthly Net Salary After Tax'
$ws.Cells.Item(55, 3).Value = 287472
$ws.Cells.Item(55, 4).Value = 287472
$ws.Cells.Item(55, 5).Value = 287472

$ws.Cells.Item(56, 2).Value = 'Mortgage Interest Rate in Percentages  Yearly for 20 Years FixedRate'
$ws.Cells.Item(56, 3).Value = 366
$ws.Cells.Item(56, 4).Value = 125
$ws.Cells.Item(56, 5).Value = 450
